$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 942.3137
$ws_ALC.Range("J17").Value = 964.449
$ws_ALC.Range("L17").Value = 2893.347
$ws_ALC.Range("N17").Value = -3229.347

# ALC row 129
$ws_ALC.Range("H129").Value = 1049.7894
$ws_ALC.Range("I129").Value = 431.42856
$ws_ALC.Range("J129").Value = 1410.5
$ws_ALC.Range("K129").Value = 1294.28568
$ws_ALC.Range("L129").Value = 4231.5
$ws_ALC.Range("M129").Value = 3705.71432
$ws_ALC.Range("N129").Value = -14231.5

# ALC row 131
$ws_ALC.Range("H131").Value = 6898.5
$ws_ALC.Range("I131").Value = 3217.9
$ws_ALC.Range("J131").Value = 16100
$ws_ALC.Range("K131").Value = 9653.700000000001
$ws_ALC.Range("L131").Value = 48300
$ws_ALC.Range("M131").Value = -4613.700000000001
$ws_ALC.Range("N131").Value = -58380

# ALC row 132
$ws_ALC.Range("H132").Value = 226128.55
$ws_ALC.Range("I132").Value = 270782.28
$ws_ALC.Range("K132").Value = 812346.8400000001
$ws_ALC.Range("M132").Value = -809816.8400000001

# ALC row 133
$ws_ALC.Range("H133").Value = 13265
$ws_ALC.Range("J133").Value = 13265
$ws_ALC.Range("L133").Value = 13265
$ws_ALC.Range("N133").Value = -23385

# ALC row 135
$ws_ALC.Range("H135").Value = 3872.6667
$ws_ALC.Range("I135").Value = 4550
$ws_ALC.Range("J135").Value = 2518
$ws_ALC.Range("K135").Value = 40950
$ws_ALC.Range("L135").Value = 22662
$ws_ALC.Range("M135").Value = -38415
$ws_ALC.Range("N135").Value = -27732

# ARM row 32
$ws_ARM.Range("H32").Value = 3866.1592
$ws_ARM.Range("I32").Value = 2739.1943
$ws_ARM.Range("K32").Value = 2739.1943
$ws_ARM.Range("M32").Value = -2452.1943

# ARM row 61
$ws_ARM.Range("H61").Value = 3203.7058
$ws_ARM.Range("I61").Value = 1962.1
$ws_ARM.Range("J61").Value = 4977.4287
$ws_ARM.Range("K61").Value = 1962.1
$ws_ARM.Range("L61").Value = 4977.4287
$ws_ARM.Range("M61").Value = -1750.1
$ws_ARM.Range("N61").Value = -5401.4287

# ARM row 102
$ws_ARM.Range("H102").Value = 2252.75
$ws_ARM.Range("I102").Value = 2000
$ws_ARM.Range("J102").Value = 2337
$ws_ARM.Range("K102").Value = 2000
$ws_ARM.Range("L102").Value = 2337
$ws_ARM.Range("M102").Value = -378
$ws_ARM.Range("N102").Value = -5581

# ARM row 122
$ws_ARM.Range("H122").Value = 1980
$ws_ARM.Range("I122").Value = 0
$ws_ARM.Range("J122").Value = 1980
$ws_ARM.Range("K122").Value = 0
$ws_ARM.Range("L122").Value = 5940
$ws_ARM.Range("M122").ClearContents()
$ws_ARM.Range("N122").Value = -10840

# ARM row 126
$ws_ARM.Range("H126").Value = 67667
$ws_ARM.Range("I126").Value = 67667
$ws_ARM.Range("K126").Value = 203001
$ws_ARM.Range("M126").Value = -200531

# ARM row 132
$ws_ARM.Range("H132").Value = 2196.3684
$ws_ARM.Range("I132").Value = 1701.421
$ws_ARM.Range("J132").Value = 3186.2632
$ws_ARM.Range("K132").Value = 5104.263
$ws_ARM.Range("L132").Value = 9558.7896
$ws_ARM.Range("M132").Value = -2574.263
$ws_ARM.Range("N132").Value = -14618.7896

# ARM row 136
$ws_ARM.Range("H136").Value = 3203.7058
$ws_ARM.Range("I136").Value = 1962.1
$ws_ARM.Range("J136").Value = 4977.4287
$ws_ARM.Range("K136").Value = 5886.299999999999
$ws_ARM.Range("L136").Value = 14932.2861
$ws_ARM.Range("M136").Value = -3336.299999999999
$ws_ARM.Range("N136").Value = -20032.2861

# ARM row 139
$ws_ARM.Range("H139").Value = 0
$ws_ARM.Range("J139").Value = 0
$ws_ARM.Range("L139").Value = 0
$ws_ARM.Range("N139").ClearContents()

# BSM row 94
$ws_BSM.Range("H94").Value = 767.2
$ws_BSM.Range("I94").Value = 757.7143
$ws_BSM.Range("J94").Value = 900
$ws_BSM.Range("K94").Value = 757.7143
$ws_BSM.Range("L94").Value = 900
$ws_BSM.Range("M94").Value = -306.7143
$ws_BSM.Range("N94").Value = -1802

# BSM row 105
$ws_BSM.Range("H105").Value = 3530.15
$ws_BSM.Range("I105").Value = 3556.4375
$ws_BSM.Range("J105").Value = 3425
$ws_BSM.Range("K105").Value = 3556.4375
$ws_BSM.Range("L105").Value = 3425
$ws_BSM.Range("M105").Value = -1809.4375
$ws_BSM.Range("N105").Value = -6919

# CRP row 58
$ws_CRP.Range("H58").Value = 1332.3529
$ws_CRP.Range("I58").Value = 760.7143
$ws_CRP.Range("K58").Value = 760.7143
$ws_CRP.Range("M58").Value = -557.7143

# CRP row 136
$ws_CRP.Range("H136").Value = 1332.3529
$ws_CRP.Range("I136").Value = 760.7143
$ws_CRP.Range("K136").Value = 2282.1429
$ws_CRP.Range("M136").Value = 267.8571000000002

# CUL row 80
$ws_CUL.Range("H80").Value = 1125
$ws_CUL.Range("I80").Value = 1050
$ws_CUL.Range("K80").Value = 3150
$ws_CUL.Range("M80").Value = -2214

# CUL row 83
$ws_CUL.Range("H83").Value = 1125
$ws_CUL.Range("I83").Value = 1050
$ws_CUL.Range("K83").Value = 9450
$ws_CUL.Range("M83").Value = -4770

# GSM row 97
$ws_GSM.Range("H97").Value = 1360
$ws_GSM.Range("I97").Value = 1360
$ws_GSM.Range("J97").Value = 0
$ws_GSM.Range("K97").Value = 1360
$ws_GSM.Range("L97").Value = 0
$ws_GSM.Range("M97").Value = -864
$ws_GSM.Range("N97").ClearContents()

# GSM row 122
$ws_GSM.Range("H122").Value = 5557059
$ws_GSM.Range("I122").Value = 5557059
$ws_GSM.Range("K122").Value = 16671177
$ws_GSM.Range("M122").Value = -16668727

# GSM row 124
$ws_GSM.Range("H124").Value = 0
$ws_GSM.Range("J124").Value = 0
$ws_GSM.Range("L124").Value = 0
$ws_GSM.Range("N124").ClearContents()

# GSM row 126
$ws_GSM.Range("H126").Value = 2520
$ws_GSM.Range("I126").Value = 2325
$ws_GSM.Range("J126").Value = 2590.9092
$ws_GSM.Range("K126").Value = 6975
$ws_GSM.Range("L126").Value = 7772.7276
$ws_GSM.Range("M126").Value = -4505
$ws_GSM.Range("N126").Value = -12712.7276

# LTW row 7
$ws_LTW.Range("H7").Value = 3331.7896
$ws_LTW.Range("J7").Value = 3500
$ws_LTW.Range("L7").Value = 3500
$ws_LTW.Range("N7").Value = -3724

# LTW row 16
$ws_LTW.Range("H16").Value = 699.5625
$ws_LTW.Range("I16").Value = 371.42856
$ws_LTW.Range("J16").Value = 954.7778
$ws_LTW.Range("K16").Value = 371.42856
$ws_LTW.Range("L16").Value = 954.7778
$ws_LTW.Range("M16").Value = -201.42856
$ws_LTW.Range("N16").Value = -1294.7778

# LTW row 40
$ws_LTW.Range("H40").Value = 3500
$ws_LTW.Range("I40").Value = 0
$ws_LTW.Range("J40").Value = 3500
$ws_LTW.Range("K40").Value = 0
$ws_LTW.Range("L40").Value = 3500
$ws_LTW.Range("M40").ClearContents()
$ws_LTW.Range("N40").Value = -3772

# LTW row 122
$ws_LTW.Range("H122").Value = 3088.3076
$ws_LTW.Range("I122").Value = 2358
$ws_LTW.Range("J122").Value = 3714.2856
$ws_LTW.Range("K122").Value = 7074
$ws_LTW.Range("L122").Value = 11142.8568
$ws_LTW.Range("M122").Value = -4624
$ws_LTW.Range("N122").Value = -16042.8568

# LTW row 126
$ws_LTW.Range("H126").Value = 3331.7896
$ws_LTW.Range("J126").Value = 3500
$ws_LTW.Range("L126").Value = 10500
$ws_LTW.Range("N126").Value = -15440

# WVR row 132
$ws_WVR.Range("H132").Value = 20002446
$ws_WVR.Range("I132").Value = 35716324
$ws_WVR.Range("J132").Value = 2967.2727
$ws_WVR.Range("K132").Value = 107148972
$ws_WVR.Range("L132").Value = 8901.8181
$ws_WVR.Range("M132").Value = -107146442
$ws_WVR.Range("N132").Value = -13961.8181
